# Locate the numbered-list paragraph that currently reads:
#   "users – пользователи системы"
# and:
#   1) rewrite its own content in place to the new "Резервное копирование..."
#      paragraph (this is what keeps the existing _GoBack bookmark trailing
#      the new text, exactly like in the original paragraph), and
#   2) insert a brand new paragraph right after it (same list style) that
#      carries the original "users – пользователи системы" content (without
#      the bookmark, which stays behind with the reworded paragraph).

$d = $word.ActiveDocument

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.StartsWith("users")) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq 0) {
    Write-Host "Target paragraph not found"
} else {
    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

    # New paragraph inserted right after the target, carrying over the old
    # "users - пользователи системы" wording (no bookmark).
    $usersXml = '<w:p ' + $wNs + '>' + `
        '<w:pPr>' + `
          '<w:pStyle w:val="a3"/>' + `
          '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
          '<w:spacing w:after="0" w:line="360" w:lineRule="auto"/>' + `
          '<w:jc w:val="both"/>' + `
          '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>users</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> – пользователи системы</w:t></w:r>' + `
      '</w:p>'

    # Replacement content for the existing paragraph: the "Резервное
    # копирование..." wording, keeping the _GoBack bookmark at the end.
    $backupXml = '<w:p ' + $wNs + '>' + `
        '<w:pPr>' + `
          '<w:pStyle w:val="a3"/>' + `
          '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
          '<w:spacing w:after="0" w:line="360" w:lineRule="auto"/>' + `
          '<w:jc w:val="both"/>' + `
          '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Резервное копирование данных в учетной записи 1</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>С:Предприятие</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> необходимо для обеспечения сохранности учетной информации и бесперебойной работы системы. Оно позволяет восстановить данные при возникновении сбоев, ошибочных действиях пользователей или технических неполадках. Регулярное выполнение резервного копирования повышает надежность и безопасность ведения учета.</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
      '</w:p>'

    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphAfter()

    # Re-fetch: target paragraph keeps its index, the freshly inserted
    # (currently empty) paragraph now sits right after it.
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertXML($backupXml)

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.InsertXML($usersXml)

    Write-Host "Paragraph split complete"
}
